$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 30.01.2022 16:30"

# Row 10 (EuroOil Opuštěná) updates
# B10: new price value (numeric)
$ws.Range("B10").Value = 36.6
# C10: old price becomes the previous B10 value (numeric)
$ws.Range("C10").Value = 36.4

# D10: delta, now stored as literal text "+0.2" (not a number) with default formatting.
# Writing the string straight into .Value would get auto-coerced back into a
# number by Excel's smart-entry logic, so round-trip it through a text
# formula and paste the computed value back in as a literal.
$ws.Range("D10").Formula = '="+0.2"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)

# E10: date/time now stored as plain text (no number formatting), same trick,
# but first strip the existing date-time formatting so the pasted value
# doesn't inherit it.
$ws.Range("E10").ClearFormats()
$ws.Range("E10").Formula = '="2022-01-30 16:33:50"'
$ws.Range("E10").Copy()
$ws.Range("E10").PasteSpecial(-4163)

$excel.CutCopyMode = $false
